$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $ws.Range("D33").Formula = "=D28*(1-0.08)+sum(B32:C32)-E6-D30"
    $ws.Range("E33").Formula = "=E28*(1-0.08)+sum(B32:D32)-E6-E30"
    $ws.Range("F33").Formula = "=F28*(1-0.08)+sum(B32:E32)-E6-F30"
    $ws.Range("G33").Formula = "=G28*(1-0.08)+sum(B32:F32)-E6-G30"
}
